# Add new model "biokurs19-09" as a new row (row 10), pushing the
# existing biokurs19-10..biokurs19-17 rows down by one, and fix two
# previously mis-entered values (L16/L17, "pseudowortlesen_anzahl_richtig")
# that should equal H-I (errors were 0, so value should equal H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 ("biokurs19-10"), shifting
# all rows from 10 downward to 11 onward.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the "biokurs19-09" model data.
$ws.Range("A10").Value2 = "biokurs19-09"
$ws.Range("B10").Value2 = 137
$ws.Range("C10").Value2 = 1
$ws.Range("D10").Value2 = 0.73
$ws.Range("E10").Value2 = 0
$ws.Range("F10").Value2 = 136
$ws.Range("G10").Value2 = 83
$ws.Range("H10").Value2 = 102
$ws.Range("I10").Value2 = 1
$ws.Range("J10").Value2 = 0.98
$ws.Range("K10").Value2 = 0
$ws.Range("L10").Value2 = 101
$ws.Range("M10").Value2 = 92
$ws.Range("N10").Value2 = "A"

# Correct two data-entry errors in what are now rows 16 and 17
# (formerly rows 15 and 16): L should equal H since no errors/omissions
# were recorded (I=0, K=0).
$ws.Range("L16").Value2 = 109
$ws.Range("L17").Value2 = 85

# Match the author's final cell selection in the sheet.
$ws.Range("L18").Select()
